$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.478.22"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "1.668.57"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'238.19"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.4797"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "'0.2630"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D9").Value = "'0.06170"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("D10").Value = "1.667.34"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "'0.06988"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("D12").Value = "'14.85"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "'0.5896"
$ws.Range("E13").Value = "  -4.74%  "
$ws.Range("D14").Value = "'4.375"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "'75.00"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'0.9997"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "25.474.01"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").Value = "'0.000006750"
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").Value = "'11.45"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "1.881.21"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("D22").Value = "'4.456"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "'8.737"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").Value = "'5.285"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'137.30"
$ws.Range("E25").Value = "  +3.89%  "
$ws.Range("D26").Value = "'15.03"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'1.727"
$ws.Range("E28").Value = "  +3.82%  "
$ws.Range("D29").Value = "'104.94"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "'3.947"
$ws.Range("E30").Value = "  +5.77%  "
$ws.Range("D31").Value = "'0.07807"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'3.649"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").Value = "'0.04243"
$ws.Range("E34").Value = "  -5.05%  "
$ws.Range("D35").Value = "'2.601"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").Value = "'0.6086"
$ws.Range("E36").Value = "  +4.50%  "
$ws.Range("D37").Value = "'0.9495"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").Value = "'2.594"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "'0.8595"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'0.01482"
$ws.Range("E41").Value = "  -5.16%  "
$ws.Range("D42").Value = "'1.847"
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").Value = "'96.14"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("D44").Value = "'0.3772"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").Value = "'4.819"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").Value = "'6.199"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").Value = "'0.05245"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "'7.382"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  +0.14%  "
